# Auto-generated Excel COM-interop script applying the 09:01:18 scrape update
# to the "horarios-141" workbook (sheets: LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 09:01:18"
$ws1.Cells.Item(3, 1).Value = "Total filas: 99"

# Rows 68-104: data table refreshed with the latest scrape pass. Existing rows
# 68-96 are overwritten in place and rows 97-104 are newly appended.
$lp1912Rows = @(
    @(68, "09:01:18", "09:01", "84_COLONIA URQUIZA-ESC 49", 0, "LP1912"),
    @(69, "09:01:18", "09:01", "215A_EL PATO", 0, "LP1912"),
    @(70, "09:01:18", "09:01", "23_HERNANDEZ", 0, "LP1912"),
    @(71, "08:04:39", "09:02", "23_HERNANDEZ", 58, "LP1912"),
    @(72, "08:04:39", "09:03", "11_ETCHEVERRY", 59, "LP1912"),
    @(73, "08:41:16", "09:04", "11_ETCHEVERRY", 23, "LP1912"),
    @(74, "09:01:18", "09:06", "23_HERNANDEZ", 5, "LP1912"),
    @(75, "09:01:18", "09:10", "16_P MOR-SANTA ANA", 9, "LP1912"),
    @(76, "08:41:16", "09:11", "16_SANTA ANA", 30, "LP1912"),
    @(77, "08:41:16", "09:11", "16_P MOR-SANTA ANA", 30, "LP1912"),
    @(78, "08:04:39", "09:16", "27_EL RETIRO", 72, "LP1912"),
    @(79, "09:01:18", "09:17", "27_EL RETIRO", 16, "LP1912"),
    @(80, "09:01:18", "09:21", "26_HERNANDEZ", 20, "LP1912"),
    @(81, "08:04:39", "09:22", "17_ROMERO", 78, "LP1912"),
    @(82, "09:01:18", "09:23", "17_ROMERO", 22, "LP1912"),
    @(83, "09:01:18", "09:23", "16_SANTA ANA", 22, "LP1912"),
    @(84, "09:01:18", "09:23", "11_ETCHEVERRY", 22, "LP1912"),
    @(85, "08:41:16", "09:24", "11_ETCHEVERRY", 43, "LP1912"),
    @(86, "08:04:39", "09:25", "81_EL PELIGRO", 81, "LP1912"),
    @(87, "09:01:18", "09:32", "15_ABASTO", 31, "LP1912"),
    @(88, "09:01:18", "09:33", "10_OLMOS", 32, "LP1912"),
    @(89, "09:01:18", "09:35", "23_HERNANDEZ", 34, "LP1912"),
    @(90, "08:04:39", "09:41", "215C_EL PATO", 97, "LP1912"),
    @(91, "09:01:18", "09:42", "215C_EL PATO", 41, "LP1912"),
    @(92, "09:01:18", "09:43", "14_ABASTO", 42, "LP1912"),
    @(93, "09:01:18", "09:44", "14_ABASTO", 63, "LP1912"),
    @(94, "09:01:18", "09:47", "16_SANTA ANA", 46, "LP1912"),
    @(95, "09:01:18", "09:52", "15_ABASTO", 51, "LP1912"),
    @(96, "09:01:18", "09:53", "10_OLMOS", 52, "LP1912"),
    @(97, "09:01:18", "10:10", "16_P MOR-SANTA ANA", 69, "LP1912"),
    @(98, "08:41:16", "10:11", "16_P MOR-SANTA ANA", 90, "LP1912"),
    @(99, "09:01:18", "10:11", "10_OLMOS", 70, "LP1912"),
    @(100, "09:01:18", "10:21", "26_HERNANDEZ", 80, "LP1912"),
    @(101, "09:01:18", "10:26", "215A_EL PATO", 85, "LP1912"),
    @(102, "09:01:18", "10:42", "17_ROMERO", 101, "LP1912"),
    @(103, "09:01:18", "10:43", "14_ABASTO", 102, "LP1912"),
    @(104, "09:01:18", "10:57", "27_EL RETIRO", 116, "LP1912"),
)

foreach ($row in $lp1912Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 09:01:18"

$ws2.Cells.Item(14, 1).Value = "09:01:18"
$ws2.Cells.Item(14, 4).Value = 0

$ws2.Cells.Item(16, 1).Value = "09:01:18"
$ws2.Cells.Item(16, 4).Value = 41

$ws2.Cells.Item(17, 1).Value = "09:01:18"
$ws2.Cells.Item(17, 4).Value = 85

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 09:01:18"
$ws3.Cells.Item(3, 1).Value = "Total filas: 20"

$ws3.Cells.Item(22, 1).Value = "09:01:18"
$ws3.Cells.Item(22, 4).Value = 9

$ws3.Cells.Item(24, 1).Value = "09:01:18"
$ws3.Cells.Item(24, 4).Value = 62

# New row 25
$ws3.Cells.Item(25, 1).Value = "09:01:18"
$ws3.Cells.Item(25, 2).Value = "10:54"
$ws3.Cells.Item(25, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(25, 4).Value = 113
$ws3.Cells.Item(25, 5).Value = "L6173"

